# Add 'hole_id' index column to the 'train' sheet so cross validation can be performed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("train")

# List of hole_id values, in row order (rows 2..38 correspond to the
# 37 data rows already present on the sheet).
$holeIds = @(
    "LBU_05_09",
    "LBU_05_30",
    "LBU_05_27",
    "LBU_02_4",
    "MHZ_12_03",
    "MHZ_12_04",
    "LBU_05_11",
    "LBU_07_03",
    "LBU_05_07",
    "LBU_05_03",
    "MHZ_08_01",
    "LBU_05_26",
    "LBU_05_29",
    "MHZ_12_01",
    "LBU_05_01",
    "MHZ_08_03",
    "LBU_01_2",
    "LBU_05_24",
    "LBU_01_3",
    "LBU_05_18",
    "LBU_07_02",
    "LBU_05_13",
    "LBU_05_22",
    "LBU_05_10",
    "LBU_05_25",
    "LBU_05_16",
    "LBU_05_28",
    "MHZ_08_04",
    "MHZ_08_02",
    "LBU_05_06",
    "LBU_05_15",
    "LBU_05_17",
    "LBU_01_1",
    "LBU_05_04",
    "LBU_02_3",
    "LBU_05_20",
    "MHZ_08_05"
)

# Header cell for the new index column, matching the style already used
# by the other header cells in row 1 (bold, bordered, centered).
$ws.Range("A1").Value = "hole_id"
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Replace the numeric row index in column A (rows 2..38) with the
# corresponding hole_id text value. Keep the existing cell style.
for ($i = 0; $i -lt $holeIds.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $holeIds[$i]
}
